$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove column I ("Reason") entirely; this shifts column J ("Food Items") into column I
$ws.Range("I1").EntireColumn.Delete() | Out-Null

# Update row 2 values
$ws.Range("A2").Value = 9
$ws.Range("B2").Value = 45691.22928240741
$ws.Range("C2").Value = "Reevan"
$ws.Range("D2").Value = 66
$ws.Range("E2").Value = 429
$ws.Range("F2").Value = 8
$ws.Range("G2").Value = 2
$ws.Range("H2").Value = 10
$ws.Range("I2").Value = "Vanilla Shake (x2)"

# Update row 3 values
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 45686.22928240741
$ws.Range("C3").Value = "Karthik"
$ws.Range("D3").Value = 2
$ws.Range("E3").Value = 790
$ws.Range("F3").Value = 10
$ws.Range("G3").Value = 5
$ws.Range("H3").Value = 15
$ws.Range("I3").Value = "Vanilla Shake (x4), Mango Lassi (x2)"
